$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("editprofile")
$wsDonor = $wb.Worksheets.Item("login")

# Remove the obsolete cash-out test rows (old rows 5-12), shifting the
# remaining rows up. This leaves only the header row plus the 3
# "Max"/"stevemax" scenario rows that are being kept/retargeted.
[void]$ws.Range("A5:H12").Delete([Microsoft.Office.Interop.Excel.XlDeleteShiftDirection]::xlShiftUp)

# Drop all existing hyperlinks - they will be rebuilt below for the rows
# that remain.
[void]$ws.Hyperlinks.Delete()

# --- Row 2: Invalid email format scenario ---
$ws.Range("A2").Value = "Invalid email format"
$ws.Range("B2").Value = "Steve"
$ws.Range("C2").Value = "Max"
$ws.Range("D2").Value = "stevemax"
$ws.Range("E2").Value = "stevemax@citrof"
$ws.Range("F2").Value = 9221244785
$ws.Range("G2").Value = "'01/01/1990"
$ws.Range("H2").WrapText = $true
$ws.Range("H2").Value = "Invalid input"

# --- Row 3: Invalid mobile number format scenario ---
$ws.Range("A3").Value = "Invalid mobile number format"
$ws.Range("B3").Value = "Steve"
$ws.Range("C3").Value = "Max"
$ws.Range("D3").Value = "stevemax"
$ws.Range("E3").Value = "stevemax@citrof.com"
$ws.Range("F3").Value = 9221244785
$ws.Range("G3").Value = "'01/01/1990"
$ws.Range("H3").Value = "Please complete all necessary details"

# --- Row 4: Invalid username format scenario ---
$ws.Range("A4").Value = "Invalid username format"
$ws.Range("B4").Value = "Steve"
$ws.Range("C4").Value = "Max"
$ws.Range("D4").Value = "stevemax"
$ws.Range("E4").Value = "stevemax@citrof.com"
$ws.Range("F4").Value = 9221244785
$ws.Range("G4").Value = "'01/01/1990"
$ws.Range("H4").WrapText = $true
$ws.Range("H4").Value = "Invalid input"

# Re-create the mailto hyperlinks on the email column for the 3 remaining rows.
[void]$ws.Hyperlinks.Add($ws.Range("E2"), "mailto:stevemax@citrof")
[void]$ws.Hyperlinks.Add($ws.Range("E3"), "mailto:stevemax@citrof.com")
[void]$ws.Hyperlinks.Add($ws.Range("E4"), "mailto:stevemax@citrof.com")

# Hyperlinks.Add() re-applies its own hyperlink formatting (a duplicate style),
# so restore the plain shared hyperlink look used by the rest of the workbook.
$wsDonor.Range("C2").Copy() | Out-Null
$ws.Range("E2:E4").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats) | Out-Null

# Selection / active cell as last left by the author.
[void]$ws.Range("B12").Select()
